$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (Jan / Can / Mittelfeld / 66)
$ws.Range("A10").Value = "Jan"
$ws.Range("B10").Value = "Can"
$ws.Range("C10").Value = "Mittelfeld"
$ws.Range("D10").Value = 66

# Update selection to reflect new active cell
$ws.Range("D10").Select()
